$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 825261.4
$ws.Range("I88").Value = 895
$ws.Range("J88").Value = 952087
$ws.Range("K88").Value = 895
$ws.Range("L88").Value = 952087
$ws.Range("M88").Value = -489
$ws.Range("N88").Value = -952899
$ws.Range("H91").Value = 825261.4
$ws.Range("I91").Value = 895
$ws.Range("J91").Value = 952087
$ws.Range("K91").Value = 895
$ws.Range("L91").Value = 952087
$ws.Range("M91").Value = 509
$ws.Range("N91").Value = -954895
$ws.Range("H121").Value = 674.0714
$ws.Range("J121").Value = 674.0714
$ws.Range("L121").Value = 2022.2142
$ws.Range("N121").Value = -5516.2142
$ws.Range("H125").Value = 2630.8333
$ws.Range("I125").Value = 3983
$ws.Range("J125").Value = 1278.6666
$ws.Range("K125").Value = 35847
$ws.Range("L125").Value = 11507.9994
$ws.Range("M125").Value = -33387
$ws.Range("N125").Value = -16427.9994
$ws.Range("H127").Value = 1275.3846
$ws.Range("I127").Value = 597.7778
$ws.Range("J127").Value = 2800
$ws.Range("K127").Value = 1793.3334
$ws.Range("L127").Value = 8400
$ws.Range("M127").Value = 3166.6666
$ws.Range("N127").Value = -18320
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 4818.5386
$ws.Range("I28").Value = 4818.5386
$ws.Range("K28").Value = 4818.5386
$ws.Range("M28").Value = -4626.5386
$ws.Range("H32").Value = 3892.2876
$ws.Range("I32").Value = 3427.2239
$ws.Range("K32").Value = 3427.2239
$ws.Range("M32").Value = -3140.2239
$ws.Range("H99").Value = 4818.5386
$ws.Range("I99").Value = 4818.5386
$ws.Range("K99").Value = 4818.5386
$ws.Range("M99").Value = -1823.5386
$ws.Range("H132").Value = 1680
$ws.Range("I132").Value = 1424
$ws.Range("K132").Value = 4272
$ws.Range("M132").Value = -1742
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 13158680
$ws.Range("I94").Value = 16667349
$ws.Range("J94").Value = 1173.75
$ws.Range("K94").Value = 16667349
$ws.Range("L94").Value = 1173.75
$ws.Range("M94").Value = -16666898
$ws.Range("N94").Value = -2075.75
$ws.Range("H105").Value = 56106588
$ws.Range("I105").Value = 67327500
$ws.Range("J105").Value = 1996.6666
$ws.Range("K105").Value = 67327500
$ws.Range("L105").Value = 1996.6666
$ws.Range("M105").Value = -67325753
$ws.Range("N105").Value = -5490.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2076.2856
$ws.Range("I31").Value = 2490
$ws.Range("J31").Value = 1910.8
$ws.Range("K31").Value = 2490
$ws.Range("L31").Value = 1910.8
$ws.Range("M31").Value = -2195
$ws.Range("N31").Value = -2500.8
$ws.Range("H34").Value = 2076.2856
$ws.Range("I34").Value = 2490
$ws.Range("J34").Value = 1910.8
$ws.Range("K34").Value = 2490
$ws.Range("L34").Value = 1910.8
$ws.Range("M34").Value = -2288
$ws.Range("N34").Value = -2314.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1627.7727
$ws.Range("I5").Value = 1795.6111
$ws.Range("J5").Value = 872.5
$ws.Range("K5").Value = 5386.8333
$ws.Range("L5").Value = 2617.5
$ws.Range("M5").Value = -5274.8333
$ws.Range("N5").Value = -2841.5
$ws.Range("H9").Value = 1335.6
$ws.Range("J9").Value = 1335.6
$ws.Range("L9").Value = 4006.8
$ws.Range("N9").Value = -4454.799999999999
$ws.Range("H104").Value = 4568.8184
$ws.Range("I104").Value = 2732
$ws.Range("J104").Value = 5618.4287
$ws.Range("K104").Value = 8196
$ws.Range("L104").Value = 16855.2861
$ws.Range("M104").Value = -5575
$ws.Range("N104").Value = -22097.2861
$ws.Range("H131").Value = 15875480
$ws.Range("J131").Value = 2991.3333
$ws.Range("L131").Value = 8973.999899999999
$ws.Range("N131").Value = -19053.9999
$ws.Range("H135").Value = 1627.7727
$ws.Range("I135").Value = 1795.6111
$ws.Range("J135").Value = 872.5
$ws.Range("K135").Value = 16160.4999
$ws.Range("L135").Value = 7852.5
$ws.Range("M135").Value = -13625.4999
$ws.Range("N135").Value = -12922.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2126.6667
$ws.Range("I126").Value = 1877.1428
$ws.Range("K126").Value = 5631.428400000001
$ws.Range("M126").Value = -3161.428400000001
$ws.Range("H132").Value = 1921.6945
$ws.Range("I132").Value = 1276.7826
$ws.Range("K132").Value = 3830.3478
$ws.Range("M132").Value = -1300.3478
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 907.6667
$ws.Range("I22").Value = 717.25
$ws.Range("J22").Value = 1060
$ws.Range("K22").Value = 717.25
$ws.Range("L22").Value = 1060
$ws.Range("M22").Value = -422.25
$ws.Range("N22").Value = -1650
$ws.Range("H27").Value = 907.6667
$ws.Range("I27").Value = 717.25
$ws.Range("J27").Value = 1060
$ws.Range("K27").Value = 717.25
$ws.Range("L27").Value = 1060
$ws.Range("M27").Value = -610.25
$ws.Range("N27").Value = -1274
$ws.Range("H132").Value = 18808.225
$ws.Range("I132").Value = 1035.4166
$ws.Range("J132").Value = 47891
$ws.Range("K132").Value = 3106.2498
$ws.Range("L132").Value = 143673
$ws.Range("M132").Value = -576.2498000000001
$ws.Range("N132").Value = -148733
$ws.Range("H134").Value = 28622.223
$ws.Range("J134").Value = 28622.223
$ws.Range("L134").Value = 28622.223
$ws.Range("N134").Value = -38762.223
$ws.Range("H137").Value = 35618.43
$ws.Range("J137").Value = 35618.43
$ws.Range("L137").Value = 35618.43
$ws.Range("N137").Value = -45818.43
$ws.Range("H139").Value = 34966.668
$ws.Range("J139").Value = 34966.668
$ws.Range("L139").Value = 34966.668
$ws.Range("N139").Value = -45246.668
$ws.Range("H141").Value = 54614.445
$ws.Range("J141").Value = 54614.445
$ws.Range("L141").Value = 54614.445
$ws.Range("N141").Value = -64974.445
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 58825308
$ws.Range("I126").Value = 83335020
$ws.Range("J126").Value = 1980.6
$ws.Range("K126").Value = 250005060
$ws.Range("L126").Value = 5941.799999999999
$ws.Range("M126").Value = -250002590
$ws.Range("N126").Value = -10881.8
$ws.Range("H132").Value = 2570.6487
$ws.Range("I132").Value = 2561.697
$ws.Range("K132").Value = 7685.091
$ws.Range("M132").Value = -5155.091
